$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos table's Price (D) / Volume(1h) (E) columns.
# D13 and D22 are prefixed with a leading apostrophe so Excel keeps the
# trailing zero ("15.00" / "226.90") as text instead of normalizing them
# to the numbers 15 and 226.9.
$ws.Range('D2').Value = '38.694.72'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '2.084.94'
$ws.Range('E3').Value = '  +1.61%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '228.34'
$ws.Range('E6').Value = '  +0.86%  '
$ws.Range('D7').Value = '60.03'
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +2.20%  '
$ws.Range('D10').Value = '0.0842'
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').Value = '2.393.96'
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').Value = '''15.00'
$ws.Range('E13').Value = '  +4.01%  '
$ws.Range('D14').Value = '21.87'
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('E15').Value = '  +4.66%  '
$ws.Range('D16').Value = '5.48'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = '2.078.71'
$ws.Range('E17').Value = '  +1.61%  '
$ws.Range('D18').Value = '38.665.87'
$ws.Range('E18').Value = '  +2.46%  '
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('D20').Value = '6.01'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('E21').Value = '  +0.90%  '
$ws.Range('D22').Value = '''226.90'
$ws.Range('E22').Value = '  +1.83%  '
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').Value = '2.42'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '2.33'
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('D26').Value = '170.92'
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('D27').Value = '9.54'
$ws.Range('E27').Value = '  +2.44%  '
$ws.Range('E28').Value = '  +7.30%  '
$ws.Range('E29').Value = '  +13.00%  '
$ws.Range('D30').Value = '19.14'
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('E32').Value = '  +4.32%  '
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('E34').Value = '  +3.19%  '
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('D36').Value = '6.46'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('D38').Value = '3.54'
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('E41').Value = '  +5.71%  '
$ws.Range('D42').Value = '1.542.93'
$ws.Range('E42').Value = '  +1.12%  '
$ws.Range('D43').Value = '100.27'
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('D44').Value = '2.81'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('E45').Value = '  +3.59%  '
$ws.Range('E46').Value = '  +9.16%  '
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('E49').Value = '  +2.66%  '
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').Value = '2.282.30'
$ws.Range('E51').Value = '  +1.74%  '
